# Append 13 new flight-arrival rows (rows 114-126) to the "Main Data" sheet.
# These represent additional arrivals for "Friday, Jan 13" appended after the
# previously last row (113) of data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 114
$ws.Range("A114").Value = 113
$ws.Range("B114").Value = "Friday, Jan 13"
$ws.Range("C114").Value = "1:48 AM"
$ws.Range("D114").Value = "UNKNOWN"
$ws.Range("E114").Value = "Tenerife"
$ws.Range("F114").Value = "(TFS)"
$ws.Range("G114").Value = "Enter Air "
$ws.Range("H114").Value = "B738"
$ws.Range("I114").Value = "(SP-ESF)"
$ws.Range("J114").Value = "1:30 AM"
$ws.Range("L114").Value = "0 hours, -18 minutes"

# Row 115
$ws.Range("A115").Value = 114
$ws.Range("B115").Value = "Friday, Jan 13"
$ws.Range("C115").Value = "6:03 AM"
$ws.Range("D115").Value = "P81988"
$ws.Range("E115").Value = "Cologne"
$ws.Range("F115").Value = "(CGN)"
$ws.Range("G115").Value = "SprintAir "
$ws.Range("H115").Value = "AT73"
$ws.Range("I115").Value = "(SP-SPD)"
$ws.Range("J115").Value = "6:09 AM"
$ws.Range("L115").Value = "0 hours, 6 minutes"

# Row 116
$ws.Range("A116").Value = 115
$ws.Range("B116").Value = "Friday, Jan 13"
$ws.Range("C116").Value = "10:05 AM"
$ws.Range("D116").Value = "LH1388"
$ws.Range("E116").Value = "Frankfurt"
$ws.Range("F116").Value = "(FRA)"
$ws.Range("G116").Value = "Lufthansa "
$ws.Range("H116").Value = "CRJ9"
$ws.Range("I116").Value = "(D-ACNO)"
$ws.Range("J116").Value = "10:06 AM"
$ws.Range("L116").Value = "0 hours, 1 minutes"

# Row 117
$ws.Range("A117").Value = 116
$ws.Range("B117").Value = "Friday, Jan 13"
$ws.Range("C117").Value = "11:03 AM"
$ws.Range("D117").Value = "LPR42"
$ws.Range("E117").Value = "Warsaw"
$ws.Range("F117").Value = "(WAW)"
$ws.Range("G117").Value = "Polish Medical Air Rescue "
$ws.Range("H117").Value = "LJ75"
$ws.Range("I117").Value = "(SP-MXS)"
$ws.Range("J117").Value = "10:42 AM"
$ws.Range("L117").Value = "0 hours, -21 minutes"

# Row 118
$ws.Range("A118").Value = 117
$ws.Range("B118").Value = "Friday, Jan 13"
$ws.Range("C118").Value = "11:05 AM"
$ws.Range("D118").Value = "FR7941"
$ws.Range("E118").Value = "Edinburgh"
$ws.Range("F118").Value = "(EDI)"
$ws.Range("G118").Value = "Ryanair "
$ws.Range("H118").Value = "B738"
$ws.Range("I118").Value = "(SP-RSM)"
$ws.Range("J118").Value = "10:55 AM"
$ws.Range("L118").Value = "0 hours, -10 minutes"

# Row 119
$ws.Range("A119").Value = 118
$ws.Range("B119").Value = "Friday, Jan 13"
$ws.Range("C119").Value = "11:40 AM"
$ws.Range("D119").Value = "FR7943"
$ws.Range("E119").Value = "Manchester"
$ws.Range("F119").Value = "(MAN)"
$ws.Range("G119").Value = "Ryanair "
$ws.Range("H119").Value = "B738"
$ws.Range("I119").Value = "(SP-RSX)"
$ws.Range("J119").Value = "11:21 AM"
$ws.Range("L119").Value = "0 hours, -19 minutes"

# Row 120
$ws.Range("A120").Value = 119
$ws.Range("B120").Value = "Friday, Jan 13"
$ws.Range("C120").Value = "11:40 AM"
$ws.Range("D120").Value = "W61926"
$ws.Range("E120").Value = "Eindhoven"
$ws.Range("F120").Value = "(EIN)"
$ws.Range("G120").Value = "Wizz Air "
$ws.Range("H120").Value = "A320"
$ws.Range("I120").Value = "(HA-LYE)"
$ws.Range("J120").Value = "11:26 AM"
$ws.Range("L120").Value = "0 hours, -14 minutes"

# Row 121
$ws.Range("A121").Value = 120
$ws.Range("B121").Value = "Friday, Jan 13"
$ws.Range("C121").Value = "12:05 PM"
$ws.Range("D121").Value = "FR1750"
$ws.Range("E121").Value = "London"
$ws.Range("F121").Value = "(STN)"
$ws.Range("G121").Value = "Ryanair "
$ws.Range("H121").Value = "B738"
$ws.Range("I121").Value = "(SP-RKR)"
$ws.Range("J121").Value = "11:38 AM"
$ws.Range("L121").Value = "0 hours, -27 minutes"

# Row 122
$ws.Range("A122").Value = 121
$ws.Range("B122").Value = "Friday, Jan 13"
$ws.Range("C122").Value = "12:20 PM"
$ws.Range("D122").Value = "LH1636"
$ws.Range("E122").Value = "Munich"
$ws.Range("F122").Value = "(MUC)"
$ws.Range("G122").Value = "Lufthansa "
$ws.Range("H122").Value = "CRJ9"
$ws.Range("I122").Value = "(D-ACNM)"
$ws.Range("J122").Value = "12:14 PM"
$ws.Range("L122").Value = "0 hours, -6 minutes"

# Row 123
$ws.Range("A123").Value = 122
$ws.Range("B123").Value = "Friday, Jan 13"
$ws.Range("C123").Value = "12:40 PM"
$ws.Range("D123").Value = "W61922"
$ws.Range("E123").Value = "Paris"
$ws.Range("F123").Value = "(BVA)"
$ws.Range("G123").Value = "Wizz Air "
$ws.Range("H123").Value = "A320"
$ws.Range("I123").Value = "(HA-LWV)"
$ws.Range("J123").Value = "12:28 PM"
$ws.Range("L123").Value = "0 hours, -12 minutes"

# Row 124
$ws.Range("A124").Value = 123
$ws.Range("B124").Value = "Friday, Jan 13"
$ws.Range("C124").Value = "1:55 PM"
$ws.Range("D124").Value = "SK1755"
$ws.Range("E124").Value = "Copenhagen"
$ws.Range("F124").Value = "(CPH)"
$ws.Range("G124").Value = "SAS "
$ws.Range("H124").Value = "CRJ9"
$ws.Range("I124").Value = "(EI-FPU)"
$ws.Range("J124").Value = "1:47 PM"
$ws.Range("L124").Value = "0 hours, -8 minutes"

# Row 125
$ws.Range("A125").Value = 124
$ws.Range("B125").Value = "Friday, Jan 13"
$ws.Range("C125").Value = "2:15 PM"
$ws.Range("D125").Value = "LO3943"
$ws.Range("E125").Value = "Warsaw"
$ws.Range("F125").Value = "(WAW)"
$ws.Range("G125").Value = "LOT "
$ws.Range("H125").Value = "E170"
$ws.Range("I125").Value = "(SP-LDF)"
$ws.Range("J125").Value = "2:03 PM"
$ws.Range("L125").Value = "0 hours, -12 minutes"

# Row 126
$ws.Range("A126").Value = 125
$ws.Range("B126").Value = "Friday, Jan 13"
$ws.Range("C126").Value = "2:35 PM"
$ws.Range("D126").Value = "KL1273"
$ws.Range("E126").Value = "Amsterdam"
$ws.Range("F126").Value = "(AMS)"
$ws.Range("G126").Value = "KLM "
$ws.Range("H126").Value = "E75L"
$ws.Range("I126").Value = "(PH-EXS)"
$ws.Range("J126").Value = "2:23 PM"
$ws.Range("L126").Value = "0 hours, -12 minutes"
